# Commit: "Add key viewer option for KeyLimiter"
#
# 1. Update Korean translations on the JudgmentVisuals sheet (NAME + DESCRIPTION rows).
# 2. Add a new "key viewer" option block (KEY/ENGLISH only, no KO/ES translations yet)
#    to the KeyLimiter sheet, right after the existing CHANGE_KEYS row.

$wb = $excel.ActiveWorkbook

# --- 1. JudgmentVisuals: refresh Korean strings -----------------------------
$wsJV = $wb.Worksheets.Item("JudgmentVisuals")
$wsJV.Range("C2").Value = "판정 비주얼"
$wsJV.Range("C3").Value = "판정을 더 정확하고 간단하게 보이도록 변경합니다."

# --- 2. KeyLimiter: append the new key-viewer rows --------------------------
$wsKL = $wb.Worksheets.Item("KeyLimiter")

$newRows = @(
    @("SHOW_KEY_VIEWER", "Show key viewer for registered keys"),
    @("KEY_VIEWER_SIZE", "Size:"),
    @("KEY_VIEWER_X_POS", "X Position:"),
    @("KEY_VIEWER_Y_POS", "Y Position:"),
    @("PRESSED_OUTLINE_COLOR", "Pressed outline color:"),
    @("RELEASED_OUTLINE_COLOR", "Released outline color:"),
    @("PRESSED_BACKGROUND_COLOR", "Pressed background color:"),
    @("RELEASED_BACKGROUND_COLOR", "Released background color:"),
    @("PRESSED_TEXT_COLOR", "Pressed text color:"),
    @("RELEASED_TEXT_COLOR", "Released text color:")
)

$startRow = 8
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $wsKL.Cells.Item($r, 1).Value = $newRows[$i][0]
    $wsKL.Cells.Item($r, 2).Value = $newRows[$i][1]
}
